$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Relative" (column G) figures were being computed against the 2nd
# data source (column E); they now need to be computed as the relative
# error against the 1st data source result (column D), i.e. F/D instead
# of F/E, for the rows whose sources actually differ.
$ws.Range("G5").Value = 10.0 / 11.0    # Difference(F5=10) / 1st Source(D5=11)
$ws.Range("G8").Value = 10.0 / 12.0    # Difference(F8=10) / 1st Source(D8=12)
$ws.Range("G26").Value = 73.0 / 15.0   # Difference(F26=73) / 1st Source(D26=15)

# Column G got a touch wider to accommodate the new percentages.
$ws.Columns.Item(7).ColumnWidth = 12.833333333333334
